# TestResources have been added
# Update resource path texts for the "transfers by date" and
# "transfers by account" rows, widen column A to fit the new text,
# and move the active selection up one row (A9 -> A8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = "/transactions/transfers/from/{date}"
$ws.Range("A9").Value = "/transactions/transfers/account/{fromAccount}"

$ws.Columns.Item(1).ColumnWidth = 34.285714285714285

$ws.Range("A8").Select()
